# Apply "PET calculator now can use MACA data" edit
# - Rename header labels on the "data" sheet
# - Clear out the 10%-increase/decrease precip stat columns (O:Q, U:W, AA:AC, AG:AI)
#   for most rows, and refresh row 5's values with newly computed MACA-based numbers
# - Restore the view's frozen-pane scroll position / active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# --- Header text updates (row 1) ---
$ws.Range("C1").Value = "Average temp and precip stats"
$ws.Range("O1").Value = "increase precipitation by 10% depth and temp with MACA R4.5 data"
$ws.Range("AA1").Value = "decrease precipitation 10% depth"

# --- Clear stale precipitation stat values ---
$ws.Range("O4:Q4").ClearContents()
$ws.Range("U4:W4").ClearContents()
$ws.Range("AA4:AC4").ClearContents()
$ws.Range("AG4:AI4").ClearContents()

$ws.Range("AA5:AC5").ClearContents()
$ws.Range("AG5:AI5").ClearContents()

$ws.Range("O6:Q6").ClearContents()
$ws.Range("U6:W6").ClearContents()
$ws.Range("AA6:AC6").ClearContents()
$ws.Range("AG6:AI6").ClearContents()

$ws.Range("O7:Q7").ClearContents()
$ws.Range("U7:W7").ClearContents()
$ws.Range("AA7:AC7").ClearContents()
$ws.Range("AG7:AI7").ClearContents()

$ws.Range("O8:Q8").ClearContents()
$ws.Range("U8:W8").ClearContents()
$ws.Range("AA8:AC8").ClearContents()
$ws.Range("AG8:AI8").ClearContents()

# --- Update row 5 with refreshed MACA-based stat values ---
$ws.Range("O5").Value = 1.597
$ws.Range("P5").Value = 1.642
$ws.Range("Q5").Value = 1.673

$ws.Range("U5").Value = 2.165
$ws.Range("V5").Value = 2.253
$ws.Range("W5").Value = 2.311

# --- Restore frozen pane view / active selection ---
$ws.Range("X18").Select()
